$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" header timestamp from 01:15 to 01:30
$ws.Range("F1").Value = "Last status check on: 27.01.2022 01:30"

# Row 5 (Makro): D5 delta and E5 last-check date were previously stored as
# text; convert them to proper numeric values, matching the other rows.
$ws.Range("D5").Value = -0.4

$ws.Range("E5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E5").Value = 44588.05217592593
